# "Clean up code and fix output" -- add the new "Yearly demand" sheet
# (hourly net-demand table, 3 days x 24 hours) as the last sheet in the
# workbook, matching the layout already used by the other hourly sheets
# (e.g. "DG Dispatch", "Connected Households", ...): row 1 / column A are
# bold+bordered+centered index headers, the rest are plain numbers.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab strip (sheetId 14 / rId14), not before the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Yearly demand"

function Write-RowRange($sheet, $row, $startCol, $values) {
    $n = $values.Length
    $data = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $data[0, $i] = $values[$i]
    }
    $startCell = $sheet.Cells.Item($row, $startCol)
    $endCell = $sheet.Cells.Item($row, $startCol + $n - 1)
    $sheet.Range($startCell, $endCell).Value = $data
}

# Row 1 (B1:Y1): hour-of-day headers 0..23
Write-RowRange $ws 1 2 @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23)

# Column A (A2:A4): day index 0..2
Write-RowRange $ws 2 1 @(0)
Write-RowRange $ws 3 1 @(1)
Write-RowRange $ws 4 1 @(2)

# Data rows 2-4: net demand values per hour
Write-RowRange $ws 2 2 @(-32.5,-19.5,-13,-13,-13,142.5,291.5,327,388.5,502,596,670.5,745,651,576.5,502,320.5,139,32,-117,-97.5,-78,-52,-39)
Write-RowRange $ws 3 2 @(-32.5,-19.5,-13,0,0,-19.5,0,324,486,648,729,751.5,583,567,333.5,340,243,57.99999999999999,-130,0,0,-78,0,-39)
Write-RowRange $ws 4 2 @(-32.5,-19.5,0,0,0,-19.5,0,0,81,324,567,589.5,648,567,324,162,81,0,-130,0,0,0,0,-39)

# Match the bold/centered/bordered header style used on every other hourly
# sheet (row 1 + column A), by copying the format from an existing sheet.
$styleSource = $wb.Worksheets.Item("DG Dispatch").Range("A2")
$styleSource.Copy()
$ws.Range("B1:Y1").PasteSpecial(-4122)
$ws.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Yearly demand sheet added"
